# The records in rows 27-30 were reordered (their non key-columns stayed
# the same per physical "slot", but the observation data moved between
# rows). Concretely this is equivalent to swapping the data of
# row 27 <-> row 29 and row 28 <-> row 30 for just the columns that
# actually carry per-record values (A, B, D, E, F, G, H, the
# Alder-Stadium/Aktivitet column K/M, Q, R, S). All other columns/cells
# are left completely untouched.
#
# We set the individual cells explicitly (rather than copying whole row
# ranges) so that cells we are not supposed to touch keep their exact
# original representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData($row, $data) {
    $ws.Range("A$row").Value2 = $data.A
    $ws.Range("B$row").Value2 = $data.B
    $ws.Range("D$row").Value2 = $data.D
    $ws.Range("E$row").Value2 = $data.E
    $ws.Range("F$row").Value2 = $data.F
    $ws.Range("G$row").Value2 = $data.G
    $ws.Range("H$row").Value2 = $data.H
    $ws.Range("Q$row").Value2 = $data.Q
    $ws.Range("R$row").Value2 = $data.R
    $ws.Range("S$row").Value2 = $data.S

    # Clear any existing K/M (Alder-Stadium / Aktivitet) cell, then set
    # the one that should be present (if any).
    $ws.Range("K$row").ClearContents()
    $ws.Range("M$row").ClearContents()
    if ($data.ContainsKey("K")) {
        $ws.Range("K$row").Value2 = $data.K
    }
    if ($data.ContainsKey("M")) {
        $ws.Range("M$row").Value2 = $data.M
    }
}

# New data for row 27 (previously held by row 29)
Set-RowData 27 @{
    A = 111940487
    B = 89573
    D = "NT"
    E = 5442
    F = "Tallticka"
    G = "Porodaedalea pini"
    H = "(Brot.) Murrill"
    Q = 575331
    R = 6633901
    S = 5
}

# New data for row 28 (previously held by row 30)
Set-RowData 28 @{
    A = 111940516
    B = 96735
    D = "VU"
    E = 220787
    F = "Knärot"
    G = "Goodyera repens"
    H = "(L.) R. Br."
    K = "fullt utvecklade blad"
    Q = 575281
    R = 6633794
    S = 5
}

# New data for row 29 (previously held by row 27)
Set-RowData 29 @{
    A = 111940543
    B = 5113
    D = "LC"
    E = 100526
    F = "Bronshjon"
    G = "Callidium coriaceum"
    H = "Paykull, 1800"
    M = "färska gnagspår"
    Q = 575299
    R = 6633979
    S = 5
}

# New data for row 30 (previously held by row 28)
Set-RowData 30 @{
    A = 111940495
    B = 56446
    D = "NT"
    E = 100049
    F = "Spillkråka"
    G = "Dryocopus martius"
    H = "(Linnaeus, 1758)"
    M = "äldre spår"
    Q = 575220
    R = 6633957
    S = 50
}
